# Test Data Added for Slovakia market
#
# Adds a new "Slovakia" worksheet (cloned from the existing "Germany"
# sheet, which shares the same layout/styling as every other country
# tab) and places it after "Portugal" as the last / active tab. Only
# the ticket reference cell (B4) is updated for the new market; the
# market-name cell (B2) is left as copied from Germany.

$wb = $excel.ActiveWorkbook

$germany  = $wb.Worksheets.Item("Germany")
$portugal = $wb.Worksheets.Item("Portugal")

# Clone Germany's sheet and drop the copy right after Portugal (i.e. at
# the end of the tab strip).
$germany.Copy($null, $portugal)
$newSheet = $wb.Worksheets.Item($portugal.Index + 1)
$newSheet.Name = "Slovakia"

# Fill in the Slovakia-specific ticket reference.
$newSheet.Range("B4").Value = "NGC-2930/T3220"

# The Germany sheet's own selection moved on as part of this edit too
# (set this before leaving the sheet, so it doesn't steal activation).
$germany.Activate() | Out-Null
$germany.Range("B19").Select() | Out-Null

# Make the new sheet the active tab, with its own selection. This must
# be the last activation/selection so Slovakia ends up as the selected
# tab in the saved workbook.
$newSheet.Activate() | Out-Null
$newSheet.Range("D14").Select() | Out-Null
